$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Button test completed: move the two "사운드" (sound) tasks (previously
# rows 17-18, not yet started) up to rows 14-15, mark them as completed with
# start/finish date+time, and shift the three tasks that used to be rows
# 14-16 down to rows 16-18 (renumbering the # column to stay sequential). ---

# Row 14 <- completed "map sound output" task (was row 17, now #10)
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = "사운드"
$ws.Range("D14").Value = "맵 사운드 출력"
$ws.Range("F14").Value = 44839
$ws.Range("G14").Value = 0.52083333333333337
$ws.Range("H14").Value = 44839
$ws.Range("I14").Value = 0.56041666666666667

# Row 15 <- completed "hit sound output" task (was row 18, now #11)
$ws.Range("B15").Value = 11
$ws.Range("C15").Value = "사운드"
$ws.Range("D15").Value = "타격 사운드 출력"
$ws.Range("F15").Value = 44839
$ws.Range("G15").Value = 0.52083333333333337
$ws.Range("H15").Value = 44839
$ws.Range("I15").Value = 0.56041666666666667

# Copy the date/time number formats from an already-completed row (row 4)
# onto the newly-filled F/G/H/I cells, instead of creating brand-new styles.
$ws.Range("F4").Copy()
$ws.Range("F14:F15").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G14:G15").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("H14:H15").PasteSpecial(-4122)
$ws.Range("I4").Copy()
$ws.Range("I14:I15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16 <- "stat window" task (was row 14, now #7)
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = "UI & 플레이어"
$ws.Range("D16").Value = "스텟창 제작"

# Row 17 <- "inventory & item" task (was row 15, now #8)
$ws.Range("B17").Value = 8
$ws.Range("C17").Value = "UI & 아이템"
$ws.Range("D17").Value = "인벤토리 및 아이템 제작"

# Row 18 <- "shortcut key" task (was row 16, now #9 - B18 was previously blank)
$ws.Range("B18").Value = 9
$ws.Range("C18").Value = "UI"
$ws.Range("D18").Value = "단축키 설정"

# E2 is "=TODAY()" (volatile) - left untouched so the engine's own
# recalculation advances its cached value naturally, the same way Excel
# would bump it a day when the workbook is reopened/edited later.

# Selection left where the author was last working.
$ws.Range("K13").Select()
